# Adapt column header formatting to respective input file names.
# - Rename "<name>_old" headers to "<name>_FV2404"
# - Rename "<name>_new" headers to "<name>_FV2410"
# - Turn the data range into an Excel Table (ListObject)
# - Freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-J (1-10) are the "old" / FV2404 columns, K (11) is "diff",
# columns L-U (12-21) are the "new" / FV2410 columns.
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value2 = ($cell.Value2 -replace '_old$', '_FV2404')
}

for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value2 = ($cell.Value2 -replace '_new$', '_FV2410')
}

# Convert the used range into a table, matching the sheet's dimensions (A1:U76).
$dataRange = $ws.Range("A1:U76")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# Freeze the header row (row 1).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
